$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the formatting of the
# existing header row (bold + border + centered), by copying H1's format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new columns I (I0) and J (IF), rows 2-76
$csv = @"
    8,8
    9,9
    7,7
    8,8
    9,9
    8,8
    9,9
    9,9
    8,8
    9,9
    7,7
    4,5
    7,7
    5,6
    7,7
    8,8
    8,8
    8,8
    7,7
    8,8
    7,7
    7,7
    7,7
    7,7
    8,8
    7,7
    5,6
    8,8
    8,9
    7,7
    9,9
    8,8
    8,8
    7,7
    8,8
    9,9
    9,9
    9,9
    6,7
    6,6
    8,9
    7,8
    6,6
    4,4
    7,7
    8,8
    6,6
    6,7
    7,8
    8,8
    12,12
    6,7
    4,5
    7,7
    5,5
    6,6
    9,9
    7,8
    8,8
    8,8
    7,7
    8,9
    8,8
    7,8
    8,9
    7,7
    7,7
    7,7
    5,5
    6,6
    5,5
    4,4
    4,4
    3,3
    6,6
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim() -ne "" }

$data = New-Object 'object[,]' $lines.Count, 2
for ($idx = 0; $idx -lt $lines.Count; $idx++) {
    $parts = $lines[$idx].Trim().Split(",")
    $data[$idx, 0] = [int]$parts[0]
    $data[$idx, 1] = [int]$parts[1]
}

$lastRow = 1 + $lines.Count
$range = $ws.Range("I2:J$lastRow")
$range.Value = $data
